# Applies the "added additional use cases" edit to TestCases.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update test-case status cells from TODO/N/A placeholders to "done"
$ws.Range("B59").Value = "done"
$ws.Range("B60").Value = "done"
$ws.Range("B61").Value = "done"
$ws.Range("B65").Value = "done"
$ws.Range("B66").Value = "done"
$ws.Range("B82").Value = "done"
$ws.Range("B84").Value = "done"
$ws.Range("B120").Value = "done"

# Row 80 no longer has a note in column C ("x")
$ws.Range("C80").ClearContents()

# Update the saved view/selection to match the author's final cursor position
$ws.Application.ActiveWindow.ScrollRow = 84
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B111").Select()
